# Insert a new row at position 4 on the first worksheet (sheet "strategy_id-0"),
# shifting the existing rows 4-11 down to rows 5-12, then populate the new
# row 4 with the "climate_change_factor_gnrl_hydropower_availability" variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 4 (elasticity_gnrl_rate_occupancy_to_gdppc),
# pushing it and everything below down by one row.
$ws.Rows("4:4").Insert()

# --- Populate the newly-inserted row 4 ---
# (columns C:G stay blank for this row, same as the surrounding data rows)
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"

$cols = @("H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS")
$vals = @(1,0.5,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $addr = "$($cols[$i])4"
    $ws.Range($addr).Value = $vals[$i]
}
